# Add more subjects to the roster on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ordered original rows plus the newly added students (name, gender, age).
$data = @(
  @("董明泽","男",7),
  @("黄毅航","男",5),
  @("刘洋","女",12),
  @("任庆阳","男",9),
  @("耿浩翔","男",6),
  @("张晓瑜","女",9),
  @("申钰嘉","男",7),
  @("郑婉怡","女",12),
  @("薛佳菲","女",12),
  @("梁文祺","男",5),
  @("陈昱","男",8),
  @("刘思潼","女",12),
  @("牛艺惠","女",6),
  @("杨舜雅","女",9),
  @("李萱","女",11),
  @("胡茗恺","男",11),
  @("孙悦媛","女",9),
  @("曾奕然","女",5),
  @("常克帅","男",15),
  @("李逸涵","男",12)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
}

# Page setup tweaks that accompanied the data entry.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author last clicked.
$ws.Range("D8").Select()
